$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Reorder "Recorded By" email list for row 12 and row 34
$ws.Range("G12").Value = "Safa.hany@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("G34").Value = "Safa.hany@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"

# Update PHYSIOLOGY session dates (rows 19-23) and session numbers (rows 21-23)
# Force text formatting so values stay as plain strings (matching source inlineStr cells)
$ws.Range("D21:D23").NumberFormat = "@"
$ws.Range("E19:E23").NumberFormat = "@"

$ws.Range("E19").Value = "13/10/2025"
$ws.Range("E20").Value = "21/10/2025"
$ws.Range("D21").Value = "3"
$ws.Range("E21").Value = "04/11/2025"
$ws.Range("D22").Value = "4"
$ws.Range("E22").Value = "06/11/2025"
$ws.Range("D23").Value = "5"
